$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A37").Value = "BRICS Initiatives for Critical Agrarian Studies (BICAS), MOSAIC Research Project, Land Deal Politics Initiative (LDPI), RCSD Chiang Mai University, Transnational Institute"
$ws.Range("B37").Value = "BRICS Initiatives in Critical Agrarian Studies;International Institute of Social Studies;Land Deal Politics Initiative;Regional Center for Social Science and Sustainable Development - Chiang Mai University;Transnational Institute"

$ws.Range("B37").Select()
